$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Data: 4 rows x 12 columns of generated YCbCr / GLCM probability values.
# ---------------------------------------------------------------------------
$values = New-Object 'object[,]' 4,12

$values[0,0]  = 0.10630677940686548
$values[0,1]  = 0.86609553710305054
$values[0,2]  = 0.36050512972792459
$values[0,3]  = 0.94686908136939518
$values[0,4]  = 0.022711554485866002
$values[0,5]  = 0.89349970342383622
$values[0,6]  = 0.76968848253133815
$values[0,7]  = 0.98864422275706709
$values[0,8]  = 0.051993135255543665
$values[0,9]  = 0.82780601953128785
$values[0,10] = 0.65043311657995173
$values[0,11] = 0.9740034323722282

$values[1,0]  = 0.09676101542419635
$values[1,1]  = 0.80214125406463943
$values[1,2]  = 0.5681229761493265
$values[1,3]  = 0.95239844174284183
$values[1,4]  = 0.0056349936620923163
$values[1,5]  = 0.83007870458569999
$values[1,6]  = 0.94366385423174304
$values[1,7]  = 0.99718250316895374
$values[1,8]  = 0.0016861844606042258
$values[1,9]  = 0.77880773414533544
$values[1,10] = 0.98748851231373314
$values[1,11] = 0.9991569077696979

$values[2,0]  = 0.12288506589364483
$values[2,1]  = 0.83987364620474714
$values[2,2]  = 0.48661755372375776
$values[2,3]  = 0.93968701951273947
$values[2,4]  = 0.025190063045391266
$values[2,5]  = 0.88304733015417569
$values[2,6]  = 0.76521315720862093
$values[2,7]  = 0.9874049684773043
$values[2,8]  = 0.00037410975217910244
$values[2,9]  = 0.21846790332599825
$values[2,10] = 0.9990935335055704
$values[2,11] = 0.99981294512391072

$values[3,0]  = 0.17581836609771587
$values[3,1]  = 0.74157535184078982
$values[3,2]  = 0.4070743092065055
$values[3,3]  = 0.91368114986194471
$values[3,4]  = 0.054561088626465429
$values[3,5]  = 0.87552188922269891
$values[3,6]  = 0.51626103829954073
$values[3,7]  = 0.97271945568676732
$values[3,8]  = 0.00010878225993910345
$values[3,9]  = 0.78361310773859749
$values[3,10] = 0.99933503133044299
$values[3,11] = 0.99994560887003037

$ws.Range("A1:L4").Value = $values

# ---------------------------------------------------------------------------
# Column widths (custom, in characters -> OOXML width units).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 12.83   # -> 13.71 chars wide
$ws.Columns.Item(2).ColumnWidth  = 11.83   # -> 12.71 chars wide
$ws.Columns.Item(3).ColumnWidth  = 11.83
$ws.Columns.Item(4).ColumnWidth  = 11.83
$ws.Columns.Item(5).ColumnWidth  = 13.83   # -> 14.71 chars wide
$ws.Columns.Item(6).ColumnWidth  = 11.83
$ws.Columns.Item(7).ColumnWidth  = 11.83
$ws.Columns.Item(8).ColumnWidth  = 11.83
$ws.Columns.Item(9).ColumnWidth  = 14.83   # -> 15.71 chars wide
$ws.Columns.Item(10).ColumnWidth = 11.83
$ws.Columns.Item(11).ColumnWidth = 11.83
$ws.Columns.Item(12).ColumnWidth = 11.83

# ---------------------------------------------------------------------------
# Style table: register a text format (numFmtId 49) and a date/time format
# (numFmtId 22) in the workbook's style table via a scratch cell, then
# clear the scratch cell so it doesn't affect the used range / data.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z100:Z101")
$ws.Range("Z100").NumberFormat = "@"
$ws.Range("Z101").NumberFormat = "m/d/yy h:mm"
$scratch.ClearFormats()
$scratch.ClearContents()

# ---------------------------------------------------------------------------
# Force a full recalculation on load (mirrors calcPr/fullCalcOnLoad="true").
# ---------------------------------------------------------------------------
$wb.RefreshAll()
